# ============================================================================
# Edit: add six new diary days (2025-06-05 .. 2025-06-09) as new date columns
# BQ:BU, populate "식단"/"대화" activity markers for those columns across the
# existing participants, and append a new participant row (구도원) with her
# own activity history.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Extend the header row (row 1) with the five new date columns.
#    Copy the formatting of the last existing date header (BP1) so the
#    new headers keep the bold / bordered / centered header style, then
#    write the date text with a leading apostrophe so Excel stores it as
#    literal text instead of auto-converting it to a date serial value.
# ----------------------------------------------------------------------
$ws.Range("BP1").Copy() | Out-Null
$ws.Range("BQ1:BU1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$newHeaderDates = @(
    ,@(69, "2025-06-05")
    ,@(70, "2025-06-06")
    ,@(71, "2025-06-07")
    ,@(72, "2025-06-08")
    ,@(73, "2025-06-09")
)
foreach ($h in $newHeaderDates) {
    $col = $h[0]
    $val = $h[1]
    $ws.Cells.Item(1, $col).Value = "'" + $val
}

# ----------------------------------------------------------------------
# 2) Populate the "식단"/"대화 (..kg)" activity cells that fall in the
#    new BQ:BU columns for existing rows (3-48), the handful of BP-column
#    entries that were previously blank, and the new row 49 (구도원),
#    including her name/handle in columns A/B.
#    Triples are (row, column, value) using 1-based indices.
# ----------------------------------------------------------------------
$bodyCells = @(
    ,@(3, 71, "식단")
    ,@(3, 72, "식단")
    ,@(6, 72, "대화 (111.0kg)")
    ,@(10, 68, "대화")
    ,@(10, 70, "식단")
    ,@(10, 71, "대화")
    ,@(10, 72, "식단")
    ,@(11, 70, "식단")
    ,@(11, 71, "식단")
    ,@(11, 72, "식단")
    ,@(12, 69, "식단")
    ,@(12, 70, "식단")
    ,@(12, 72, "식단")
    ,@(13, 72, "식단")
    ,@(13, 73, "식단")
    ,@(17, 69, "식단")
    ,@(17, 71, "식단")
    ,@(19, 69, "식단")
    ,@(19, 70, "식단")
    ,@(19, 71, "식단")
    ,@(19, 72, "식단")
    ,@(22, 72, "식단")
    ,@(24, 70, "식단")
    ,@(24, 72, "대화 (101.0kg)")
    ,@(25, 72, "식단")
    ,@(26, 68, "식단")
    ,@(26, 69, "식단")
    ,@(26, 70, "식단")
    ,@(26, 72, "식단")
    ,@(30, 70, "대화")
    ,@(30, 72, "대화 (120.0kg)")
    ,@(31, 68, "대화")
    ,@(31, 69, "식단")
    ,@(31, 70, "대화")
    ,@(31, 71, "대화")
    ,@(31, 72, "식단")
    ,@(32, 68, "식단")
    ,@(32, 70, "식단")
    ,@(35, 68, "식단")
    ,@(35, 70, "식단")
    ,@(35, 71, "식단")
    ,@(37, 70, "대화")
    ,@(37, 71, "식단")
    ,@(37, 72, "대화 (79.0kg)")
    ,@(41, 70, "식단")
    ,@(42, 68, "대화")
    ,@(42, 69, "식단")
    ,@(42, 70, "대화")
    ,@(42, 71, "식단")
    ,@(42, 72, "대화 (97.0kg)")
    ,@(42, 73, "식단")
    ,@(43, 68, "식단")
    ,@(43, 69, "식단")
    ,@(43, 70, "식단")
    ,@(43, 71, "식단")
    ,@(43, 72, "식단")
    ,@(43, 73, "식단")
    ,@(44, 70, "식단")
    ,@(45, 68, "식단")
    ,@(45, 72, "식단")
    ,@(49, 1, "구도원")
    ,@(49, 2, "gudowon01")
    ,@(49, 62, "식단")
    ,@(49, 63, "식단")
    ,@(49, 64, "식단")
    ,@(49, 65, "식단")
    ,@(49, 66, "식단")
    ,@(49, 67, "식단")
    ,@(49, 68, "식단")
    ,@(49, 70, "식단")
    ,@(49, 71, "식단")
    ,@(49, 72, "식단")
)

foreach ($cell in $bodyCells) {
    $r = $cell[0]
    $c = $cell[1]
    $v = $cell[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Host "Applied diary update: new date columns BQ:BU, activity markers, and new participant row."
